# zone.xlsx — add a "trajet" column alongside the renamed "transport" column,
# and chain the running-total formulas down the bottom two rows
# ("enchainement des methodes des dataframes").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
# B1 was "indemnité" -> "transport"; new C1 header "trajet"
$ws.Range("B1").Value = "transport"
$ws.Range("C1").Value = "trajet"

# --- Column B: updated literal values (rows 2-6) ---------------------------
$ws.Range("B2").Value = 2.64
$ws.Range("B3").Value = 6.26
$ws.Range("B4").Value = 9.69
$ws.Range("B5").Value = 12.73
$ws.Range("B6").Value = 16.28

# --- Column C: brand new literal values (rows 2-6) --------------------------
$ws.Range("C2").Value = 1.55
$ws.Range("C3").Value = 3.11
$ws.Range("C4").Value = 4.66
$ws.Range("C5").Value = 6.22
$ws.Range("C6").Value = 7.78

# --- Rows 7-8: running totals via formulas instead of literals -------------
$ws.Range("B7").Formula = "=B6+B2"
$ws.Range("C7").Formula = "=C6+C2"
$ws.Range("B8").Formula = "=B7+B2"
$ws.Range("C8").Formula = "=C7+C2"

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.6328125
$ws.Columns.Item(3).ColumnWidth = 9.6328125

# --- Final selection, as left by the author ---------------------------------
$ws.Range("D12").Select()
